$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "69.917.65"
Set-TextValue "E2" "  +0.02%  "
Set-TextValue "D3" "3.540.74"
Set-TextValue "E3" "  +0.15%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.17%  "
Set-TextValue "D5" "604.27"
Set-TextValue "E5" "  -1.57%  "
Set-TextValue "D6" "196.30"
Set-TextValue "E6" "  +4.84%  "
Set-TextValue "D7" "0.632"
Set-TextValue "E7" "  -0.05%  "
Set-TextValue "E8" "  -0.08%  "
Set-TextValue "D9" "0.206"
Set-TextValue "E9" "  -4.84%  "
Set-TextValue "D10" "0.652"
Set-TextValue "E10" "  -1.17%  "
Set-TextValue "D11" "53.86"
Set-TextValue "E11" "  +0.58%  "
Set-TextValue "D12" "0.0000302"
Set-TextValue "E12" "  -1.65%  "
Set-TextValue "E13" "  -1.44%  "
Set-TextValue "D14" "4.103.56"
Set-TextValue "E14" "  -0.13%  "
Set-TextValue "D15" "606.66"
Set-TextValue "E15" "  -1.77%  "
Set-TextValue "D16" "19.20"
Set-TextValue "E16" "  +0.27%  "
Set-TextValue "D17" "70.107.96"
Set-TextValue "E17" "  +0.13%  "
Set-TextValue "D18" "12.74"
Set-TextValue "E18" "  -0.78%  "
Set-TextValue "D19" "3.529.30"
Set-TextValue "E19" "  -0.29%  "
Set-TextValue "E20" "  +0.59%  "
Set-TextValue "D21" "0.996"
Set-TextValue "E21" "  -0.17%  "
Set-TextValue "D22" "18.22"
Set-TextValue "E22" "  +3.70%  "
Set-TextValue "D23" "5.25"
Set-TextValue "E23" "  +4.37%  "
Set-TextValue "D24" "102.80"
Set-TextValue "E24" "  -2.51%  "
Set-TextValue "E25" "  -1.58%  "
Set-TextValue "D26" "3.17"
Set-TextValue "E26" "  +4.34%  "
Set-TextValue "E27" "  -0.12%  "
Set-TextValue "E28" "  -3.95%  "
Set-TextValue "D29" "33.76"
Set-TextValue "E29" "  -1.41%  "
Set-TextValue "D30" "7.14"
Set-TextValue "E30" "  +1.38%  "
Set-TextValue "D31" "4.31"
Set-TextValue "E31" "  +15.59%  "
Set-TextValue "D32" "12.63"
Set-TextValue "E32" "  +1.27%  "
Set-TextValue "E33" "  -1.30%  "
Set-TextValue "D34" "63.21"
Set-TextValue "E34" "  -1.15%  "
Set-TextValue "D35" "0.0₃0855"
Set-TextValue "E35" "  +10.16%  "
Set-TextValue "D36" "3.730.75"
Set-TextValue "E36" "  +5.03%  "
Set-TextValue "D38" "3.05"
Set-TextValue "E38" "  -2.68%  "
Set-TextValue "E39" "  +1.36%  "
Set-TextValue "E40" "  -1.10%  "
Set-TextValue "D41" "36.70"
Set-TextValue "E41" "  -0.33%  "
Set-TextValue "D42" "488.96"
Set-TextValue "E42" "  -8.96%  "
Set-TextValue "E43" "  -5.12%  "
Set-TextValue "D44" "0.0460"
Set-TextValue "E44" "  -0.75%  "
Set-TextValue "E45" "  -1.73%  "
Set-TextValue "D46" "2.85"
Set-TextValue "E46" "  -3.78%  "
Set-TextValue "D47" "3.32"
Set-TextValue "E47" "  -1.02%  "
Set-TextValue "D48" "1.01"
Set-TextValue "E48" "  +0.25%  "
Set-TextValue "D49" "8.60"
Set-TextValue "E49" "  -3.95%  "
Set-TextValue "E50" "  +6.50%  "
Set-TextValue "D51" "130.76"
Set-TextValue "E51" "  -1.34%  "
